# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# matching the canonical OOXML diff for Sheets/Lamia_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 256.3889
$ws.Range("I39").Value = 179.85715
$ws.Range("K39").Value = 539.5714499999999
$ws.Range("M39").Value = -243.5714499999999
$ws.Range("H40").Value = 4090.4614
$ws.Range("I40").Value = 2619.3
$ws.Range("J40").Value = 8994.333000000001
$ws.Range("K40").Value = 2619.3
$ws.Range("L40").Value = 8994.333000000001
$ws.Range("M40").Value = -2444.3
$ws.Range("N40").Value = -9344.333000000001
$ws.Range("H43").Value = 16889.107
$ws.Range("I43").Value = 16333.277
$ws.Range("J43").Value = 17889.6
$ws.Range("K43").Value = 16333.277
$ws.Range("L43").Value = 17889.6
$ws.Range("M43").Value = -16264.277
$ws.Range("N43").Value = -18027.6
$ws.Range("H64").Value = 9325.333000000001
$ws.Range("I64").Value = 9997
$ws.Range("K64").Value = 9997
$ws.Range("M64").Value = -9749
$ws.Range("H67").Value = 9325.333000000001
$ws.Range("I67").Value = 9997
$ws.Range("K67").Value = 9997
$ws.Range("M67").Value = -9139
$ws.Range("H74").Value = 10738.25
$ws.Range("I74").Value = 8580.4
$ws.Range("K74").Value = 8580.4
$ws.Range("M74").Value = -7644.4
$ws.Range("H76").Value = 19002
$ws.Range("J76").Value = 19002
$ws.Range("L76").Value = 19002
$ws.Range("N76").Value = -19632
$ws.Range("H77").Value = 10738.25
$ws.Range("I77").Value = 8580.4
$ws.Range("K77").Value = 42902
$ws.Range("M77").Value = -38222
$ws.Range("H79").Value = 19002
$ws.Range("J79").Value = 19002
$ws.Range("L79").Value = 19002
$ws.Range("N79").Value = -21186
$ws.Range("H113").Value = 58661.95
$ws.Range("I113").Value = 168665.83
$ws.Range("J113").Value = 7890.923
$ws.Range("K113").Value = 168665.83
$ws.Range("L113").Value = 7890.923
$ws.Range("M113").Value = -165411.83
$ws.Range("N113").Value = -14398.923
$ws.Range("H116").Value = 11197.6
$ws.Range("I116").Value = 9524.25
$ws.Range("J116").Value = 13110
$ws.Range("K116").Value = 9524.25
$ws.Range("L116").Value = 13110
$ws.Range("M116").Value = -6082.25
$ws.Range("N116").Value = -19994
$ws.Range("H137").Value = 27030128
$ws.Range("I137").Value = 55557740
$ws.Range("J137").Value = 3967
$ws.Range("K137").Value = 166673220
$ws.Range("L137").Value = 11901
$ws.Range("M137").Value = -166670670
$ws.Range("N137").Value = -17001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4879.5835
$ws.Range("I45").Value = 3024.5
$ws.Range("J45").Value = 6734.6665
$ws.Range("K45").Value = 3024.5
$ws.Range("L45").Value = 6734.6665
$ws.Range("M45").Value = -2647.5
$ws.Range("N45").Value = -7488.6665
$ws.Range("H122").Value = 41668620
$ws.Range("I122").Value = 1973.762
$ws.Range("J122").Value = 333335140
$ws.Range("K122").Value = 5921.286
$ws.Range("L122").Value = 1000005420
$ws.Range("M122").Value = -3471.286
$ws.Range("N122").Value = -1000010320

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 32614.3
$ws.Range("I82").Value = 4021
$ws.Range("K82").Value = 4021
$ws.Range("M82").Value = -3638
$ws.Range("H85").Value = 32614.3
$ws.Range("I85").Value = 4021
$ws.Range("K85").Value = 4021
$ws.Range("M85").Value = -2695
$ws.Range("H86").Value = 3114.9
$ws.Range("I86").Value = 3114.9
$ws.Range("K86").Value = 3114.9
$ws.Range("M86").Value = -1991.9
$ws.Range("H89").Value = 3114.9
$ws.Range("I89").Value = 3114.9
$ws.Range("K89").Value = 15574.5
$ws.Range("M89").Value = -9958.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11798.1
$ws.Range("I86").Value = 5744.75
$ws.Range("K86").Value = 5744.75
$ws.Range("M86").Value = -4621.75
$ws.Range("H89").Value = 11798.1
$ws.Range("I89").Value = 5744.75
$ws.Range("K89").Value = 28723.75
$ws.Range("M89").Value = -23107.75
$ws.Range("H95").Value = 42507.75
$ws.Range("J95").Value = 42507.75
$ws.Range("L95").Value = 42507.75
$ws.Range("N95").Value = -47999.75
$ws.Range("H122").Value = 40824.965
$ws.Range("I122").Value = 55252.527
$ws.Range("K122").Value = 165757.581
$ws.Range("M122").Value = -163307.581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6059773
$ws.Range("I4").Value = 3809600
$ws.Range("K4").Value = 11428800
$ws.Range("M4").Value = -11428688
$ws.Range("H14").Value = 1104.8235
$ws.Range("I14").Value = 1104.8235
$ws.Range("K14").Value = 3314.4705
$ws.Range("M14").Value = -3141.4705
$ws.Range("H60").Value = 1027306.5
$ws.Range("I60").Value = 350.4
$ws.Range("J60").Value = 2738900
$ws.Range("K60").Value = 1051.2
$ws.Range("L60").Value = 8216700
$ws.Range("M60").Value = -800.1999999999998
$ws.Range("N60").Value = -8217202
$ws.Range("H81").Value = 3134.7646
$ws.Range("I81").Value = 1339.125
$ws.Range("J81").Value = 4730.8887
$ws.Range("K81").Value = 4017.375
$ws.Range("L81").Value = 14192.6661
$ws.Range("M81").Value = -2894.375
$ws.Range("N81").Value = -16438.6661
$ws.Range("H84").Value = 3134.7646
$ws.Range("I84").Value = 1339.125
$ws.Range("J84").Value = 4730.8887
$ws.Range("K84").Value = 12052.125
$ws.Range("L84").Value = 42577.99830000001
$ws.Range("M84").Value = -6436.125
$ws.Range("N84").Value = -53809.99830000001
$ws.Range("H116").Value = 5506.75
$ws.Range("I116").Value = 997.5
$ws.Range("J116").Value = 10016
$ws.Range("K116").Value = 2992.5
$ws.Range("L116").Value = 30048
$ws.Range("M116").Value = 449.5
$ws.Range("N116").Value = -36932
$ws.Range("I129").Value = 623.2222
$ws.Range("J129").Value = 7579919.5
$ws.Range("K129").Value = 1869.6666
$ws.Range("L129").Value = 22739758.5
$ws.Range("M129").Value = 3130.3334
$ws.Range("N129").Value = -22749758.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11351.375
$ws.Range("I80").Value = 8361
$ws.Range("K80").Value = 8361
$ws.Range("M80").Value = -7363
$ws.Range("H83").Value = 11351.375
$ws.Range("I83").Value = 8361
$ws.Range("K83").Value = 41805
$ws.Range("M83").Value = -36813
$ws.Range("H126").Value = 8779.223
$ws.Range("I126").Value = 4166.5
$ws.Range("K126").Value = 12499.5
$ws.Range("M126").Value = -10029.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 63429
$ws.Range("J128").Value = 63429
$ws.Range("L128").Value = 63429
$ws.Range("N128").Value = -73389

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16680667
$ws.Range("I5").Value = 22000
$ws.Range("K5").Value = 22000
$ws.Range("M5").Value = -21888
$ws.Range("H62").Value = 7633.3335
$ws.Range("I62").Value = 7450
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 7450
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -6826
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7633.3335
$ws.Range("I65").Value = 7450
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 37250
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -34130
$ws.Range("N65").Value = -46240
$ws.Range("H123").Value = 85000
$ws.Range("J123").Value = 85000
$ws.Range("L123").Value = 85000
$ws.Range("N123").Value = -94800
$ws.Range("H136").Value = 2512.1353
$ws.Range("I136").Value = 1528.6061
$ws.Range("K136").Value = 4585.8183
$ws.Range("M136").Value = -2035.8183

Write-Output "Updated 193 cells across 8 sheets"
